$wb = $excel.ActiveWorkbook

# --- Sheet "Folder Inventory": a new commit for "GitHub Copilot Innovation
# Workshop" moved from its old chronological slot to the top of the list.
$wsInv = $wb.Worksheets.Item("Folder Inventory")

# Insert a blank row at position 3 (not directly under the styled header
# row, so the new row stays plain/unstyled), then slide the current row 2
# down into it before overwriting row 2 with the new, most-recent entry.
$wsInv.Rows.Item(3).Insert()

$wsInv.Range("A3").Value = $wsInv.Range("A2").Value2
$wsInv.Range("B3").Value = $wsInv.Range("B2").Value2
$wsInv.Range("C3").Value = $wsInv.Range("C2").Value2
$wsInv.Range("D3").Value = $wsInv.Range("D2").Value2
$wsInv.Range("E3").Value = $wsInv.Range("E2").Value2

$wsInv.Range("A2").Value = "GitHub Copilot Innovation Workshop"
$wsInv.Range("B2").Value = "GitHub Copilot Innovation Workshop"
$wsInv.Range("C2").Value = "2025-06-16 14:35:58 +0530"
$wsInv.Range("D2").Value = 1
$wsInv.Range("E2").Value = "Root"

# The folder's previous entry (old row 46, now shifted to row 47 after the
# insert above) is now a stale duplicate, so remove it; this shifts the
# remaining rows back up, keeping the total row count unchanged.
$wsInv.Rows.Item(47).Delete()

# --- Sheet "Metadata": refresh the generation timestamp and workflow run.
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "2025-06-16 09:06:15 UTC"

# "Workflow Run" is stored as text (not a number) in the source file, so
# force text formatting before assigning, then drop the format again so no
# "number stored as text" styling sneaks into the saved workbook.
$wsMeta.Range("B5").NumberFormat = "@"
$wsMeta.Range("B5").Value = "24"
$wsMeta.Range("B5").ClearFormats()

# --- Sheet "Summary": reflect the newest "Last Updated" value.
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = "2025-06-16 14:35:58 +0530"
